$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.931.75'
$ws.Range('E2').Value = '  +8.27%  '
$ws.Range('D3').Value = '1.811.94'
$ws.Range('E3').Value = '  +5.04%  '
$ws.Range('D4').Value = "'0.9994"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'248.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9993"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = "'0.4949"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.43%  '
$ws.Range('D8').Value = "'0.2788"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.15%  '
$ws.Range('D9').Value = "'0.06421"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.93%  '
$ws.Range('D10').Value = '1.808.27'
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('D11').Value = "'16.80"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.73%  '
$ws.Range('D12').Value = "'0.07077"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('D13').Value = "'0.6466"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.17%  '
$ws.Range('D14').Value = "'83.81"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.90%  '
$ws.Range('D15').Value = "'4.693"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.20%  '
$ws.Range('D16').Value = '28.947.10'
$ws.Range('E16').Value = '  +9.03%  '
$ws.Range('D17').Value = "'0.9994"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = "'0.000007364"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.88%  '
$ws.Range('D19').Value = "'0.9995"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +8.20%  '
$ws.Range('D21').Value = '2.039.29'
$ws.Range('E21').Value = '  +4.96%  '
$ws.Range('D22').Value = "'4.587"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.07%  '
$ws.Range('D23').Value = "'8.823"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.20%  '
$ws.Range('D24').Value = "'5.355"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.12%  '
$ws.Range('D25').Value = "'143.20"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('D26').Value = "'129.42"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +21.57%  '
$ws.Range('D27').Value = "'16.45"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.07%  '
$ws.Range('D28').Value = "'1.891"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.50%  '
$ws.Range('D29').Value = "'1.416"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('D30').Value = "'4.144"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('D31').Value = "'0.08369"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.76%  '
$ws.Range('D32').Value = "'3.806"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.89%  '
$ws.Range('D33').Value = "'0.04949"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.83%  '
$ws.Range('D34').Value = "'1.101"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.26%  '
$ws.Range('D35').Value = "'0.6735"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.92%  '
$ws.Range('D36').Value = "'2.683"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('D37').Value = "'2.301"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +14.63%  '
$ws.Range('D38').Value = "'2.744"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.14%  '
$ws.Range('D39').Value = "'0.9574"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('E40').Value = '  +9.69%  '
$ws.Range('D41').Value = "'0.01592"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.47%  '
$ws.Range('D42').Value = "'0.9994"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = "'0.4101"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.23%  '
$ws.Range('D44').Value = "'100.27"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').Value = "'7.159"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('D47').Value = "'0.05533"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.19%  '
$ws.Range('D48').Value = "'31.82"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.77%  '
$ws.Range('D49').Value = "'8.134"
$ws.Range('D49').Style = 'Normal'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'1.312"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.65%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = "'0.3631"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.43%  '
